# "Generate Report for Handoff"
#
# The localization status report moved from "In Translation" to
# "Ready for handoff": update the per-language status cells and the
# "Latest Handoff/HO Xliff Generate" timestamps that were refreshed by the
# same run, then widen the Status columns to fit the new, longer text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-27 08:38:50"

# --- zh-cn sheet -------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-27 08:38:46"

# --- de-de sheet -------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-27 08:38:50"

# --- Widen the Status columns to fit "Ready for handoff" ---------------
# Overview: zh-cn (E) and de-de (F) status columns
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332

# zh-cn / de-de: Status column (C)
$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333332
$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333332
